# Weekly data refresh: insert 3 new price observations (one per "Variedad")
# at the top of the existing data block (row 912), pushing the rest of the
# rows (912-1005) down by 3 rows (to 915-1008).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 912, shifting existing rows 912:1005 down to 915:1008.
$ws.Rows("912:914").Insert()

# --- New row 912: Lechuga / Conconina(o) / Primera ---
$ws.Cells.Item(912, 1).Value = 4
$ws.Cells.Item(912, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(912, 3).Value = "Los Lagos"
$ws.Cells.Item(912, 4).Value = 45212
$ws.Cells.Item(912, 5).Value = 10
$ws.Cells.Item(912, 6).Value = 100112033
$ws.Cells.Item(912, 7).Value = "Lechuga"
$ws.Cells.Item(912, 8).Value = "Conconina(o)"
$ws.Cells.Item(912, 9).Value = "Primera"
$ws.Cells.Item(912, 10).Value = 200
$ws.Cells.Item(912, 11).Value = 11000
$ws.Cells.Item(912, 12).Value = 11000
$ws.Cells.Item(912, 13).Value = 11000
$ws.Cells.Item(912, 14).Value = "`$/caja 10 unidades"
$ws.Cells.Item(912, 15).Value = "Región Metropolitana"
$ws.Cells.Item(912, 16).Value = 1100
$ws.Cells.Item(912, 17).Value = 10
$ws.Cells.Item(912, 18).Value = "Hortaliza"

# --- New row 913: Lechuga / Escarola / Primera ---
$ws.Cells.Item(913, 1).Value = 4
$ws.Cells.Item(913, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(913, 3).Value = "Los Lagos"
$ws.Cells.Item(913, 4).Value = 45212
$ws.Cells.Item(913, 5).Value = 10
$ws.Cells.Item(913, 6).Value = 100112033
$ws.Cells.Item(913, 7).Value = "Lechuga"
$ws.Cells.Item(913, 8).Value = "Escarola"
$ws.Cells.Item(913, 9).Value = "Primera"
$ws.Cells.Item(913, 10).Value = 600
$ws.Cells.Item(913, 11).Value = 17000
$ws.Cells.Item(913, 12).Value = 18000
$ws.Cells.Item(913, 13).Value = 17500
$ws.Cells.Item(913, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(913, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(913, 16).Value = 1167
$ws.Cells.Item(913, 17).Value = 15
$ws.Cells.Item(913, 18).Value = "Hortaliza"

# --- New row 914: Lechuga / Marina / Primera ---
$ws.Cells.Item(914, 1).Value = 4
$ws.Cells.Item(914, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(914, 3).Value = "Los Lagos"
$ws.Cells.Item(914, 4).Value = 45212
$ws.Cells.Item(914, 5).Value = 10
$ws.Cells.Item(914, 6).Value = 100112033
$ws.Cells.Item(914, 7).Value = "Lechuga"
$ws.Cells.Item(914, 8).Value = "Marina"
$ws.Cells.Item(914, 9).Value = "Primera"
$ws.Cells.Item(914, 10).Value = 300
$ws.Cells.Item(914, 11).Value = 12000
$ws.Cells.Item(914, 12).Value = 12000
$ws.Cells.Item(914, 13).Value = 12000
$ws.Cells.Item(914, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(914, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(914, 16).Value = 800
$ws.Cells.Item(914, 17).Value = 15
$ws.Cells.Item(914, 18).Value = "Hortaliza"
